$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (and a few re-ordered rows).
# A leading apostrophe forces Excel to store the value as literal text,
# matching the original inline-string (non-numeric) cell content exactly.

# Row 2
$ws.Range("D2").Value = "'62.570.18"
$ws.Range("E2").Value = "'  -1.73%  "

# Row 3
$ws.Range("D3").Value = "'3.185.79"
$ws.Range("E3").Value = "'  -3.58%  "

# Row 4
$ws.Range("E4").Value = "'  +0.15%  "

# Row 5
$ws.Range("D5").Value = "'587.72"
$ws.Range("E5").Value = "'  -2.68%  "

# Row 6
$ws.Range("D6").Value = "'135.95"
$ws.Range("E6").Value = "'  -4.46%  "

# Row 7
$ws.Range("E7").Value = "'  +0.02%  "

# Row 8
$ws.Range("D8").Value = "'3.181.86"
$ws.Range("E8").Value = "'  -3.66%  "

# Row 9
$ws.Range("D9").Value = "'0.507"
$ws.Range("E9").Value = "'  -2.48%  "

# Row 10
$ws.Range("E10").Value = "'  -5.04%  "

# Row 11
$ws.Range("D11").Value = "'5.26"
$ws.Range("E11").Value = "'  -4.19%  "

# Row 12
$ws.Range("E12").Value = "'  -3.64%  "

# Row 13
$ws.Range("D13").Value = "'0.0000235"
$ws.Range("E13").Value = "'  -5.00%  "

# Row 14
$ws.Range("D14").Value = "'33.34"
$ws.Range("E14").Value = "'  -3.87%  "

# Row 15
$ws.Range("D15").Value = "'3.714.45"
$ws.Range("E15").Value = "'  -3.43%  "

# Row 16
$ws.Range("E16").Value = "'  -2.03%  "

# Row 17
$ws.Range("D17").Value = "'3.192.13"
$ws.Range("E17").Value = "'  -3.19%  "

# Row 18
$ws.Range("D18").Value = "'62.599.39"
$ws.Range("E18").Value = "'  -1.82%  "

# Row 19
$ws.Range("D19").Value = "'6.52"
$ws.Range("E19").Value = "'  -5.20%  "

# Row 20
$ws.Range("D20").Value = "'456.32"
$ws.Range("E20").Value = "'  -5.09%  "

# Row 21
$ws.Range("D21").Value = "'13.92"
$ws.Range("E21").Value = "'  -1.55%  "

# Row 22
$ws.Range("D22").Value = "'0.703"
$ws.Range("E22").Value = "'  -4.06%  "

# Row 23
$ws.Range("D23").Value = "'7.62"
$ws.Range("E23").Value = "'  -5.13%  "

# Row 24
$ws.Range("B24").Value = "'Litecoin"
$ws.Range("C24").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'83.63"
$ws.Range("E24").Value = "'  -1.33%  "

# Row 25
$ws.Range("B25").Value = "'InternetComputer(DFINITY)"
$ws.Range("C25").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "'13.22"
$ws.Range("E25").Value = "'  -2.12%  "

# Row 26
$ws.Range("E26").Value = "'  +0.00%  "

# Row 27
$ws.Range("D27").Value = "'2.69"
$ws.Range("E27").Value = "'  -2.72%  "

# Row 28
$ws.Range("E28").Value = "'  +0.10%  "

# Row 29
$ws.Range("B29").Value = "'RenderToken"
$ws.Range("C29").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'7.78"
$ws.Range("E29").Value = "'  -4.48%  "

# Row 30
$ws.Range("B30").Value = "'NEARProtocol"
$ws.Range("C30").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'6.82"
$ws.Range("E30").Value = "'  -6.56%  "

# Row 31
$ws.Range("E31").Value = "'  -7.04%  "

# Row 32
$ws.Range("D32").Value = "'27.34"
$ws.Range("E32").Value = "'  -6.48%  "

# Row 33
$ws.Range("D33").Value = "'0.103"
$ws.Range("E33").Value = "'  -2.09%  "

# Row 34
$ws.Range("D34").Value = "'2.37"
$ws.Range("E34").Value = "'  -6.22%  "

# Row 35
$ws.Range("E35").Value = "'  -5.25%  "

# Row 36
$ws.Range("D36").Value = "'5.89"
$ws.Range("E36").Value = "'  -1.44%  "

# Row 37
$ws.Range("D37").Value = "'51.09"
$ws.Range("E37").Value = "'  -3.50%  "

# Row 38
$ws.Range("D38").Value = "'0.0₃0699"
$ws.Range("E38").Value = "'  -6.40%  "

# Row 39
$ws.Range("D39").Value = "'0.0385"
$ws.Range("E39").Value = "'  -3.97%  "

# Row 40
$ws.Range("E40").Value = "'  -0.51%  "

# Row 41
$ws.Range("B41").Value = "'Cosmos"
$ws.Range("C41").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").Value = "'7.99"
$ws.Range("E41").Value = "'  -4.71%  "

# Row 42
$ws.Range("B42").Value = "'Maker"
$ws.Range("C42").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "'2.834.87"
$ws.Range("E42").Value = "'  -7.37%  "

# Row 43
$ws.Range("B43").Value = "'Kaspa"
$ws.Range("C43").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.111"
$ws.Range("E43").Value = "'  -0.15%  "

# Row 44
$ws.Range("B44").Value = "'Bittensor"
$ws.Range("C44").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").Value = "'390.75"
$ws.Range("E44").Value = "'  -8.30%  "

# Row 45
$ws.Range("D45").Value = "'36.26"
$ws.Range("E45").Value = "'  +3.69%  "

# Row 46
$ws.Range("D46").Value = "'0.249"
$ws.Range("E46").Value = "'  -6.40%  "

# Row 47
$ws.Range("B47").Value = "'USDe"
$ws.Range("C47").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "'0.999"
$ws.Range("E47").Value = "'  -0.04%  "

# Row 48
$ws.Range("B48").Value = "'Fetch.AI"
$ws.Range("C48").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").Value = "'2.13"
$ws.Range("E48").Value = "'  -3.29%  "

# Row 49
$ws.Range("D49").Value = "'124.97"
$ws.Range("E49").Value = "'  +0.65%  "

# Row 50
$ws.Range("D50").Value = "'25.46"
$ws.Range("E50").Value = "'  -3.25%  "

# Row 51
$ws.Range("E51").Value = "'  -3.73%  "
